$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.73
$ws.Range("H2").Value = 4.1
$ws.Range("I2").Value = 4.2
$ws.Range("Q2").Value = 1.53
$ws.Range("R2").Value = 2.5
$ws.Range("S2").Value = 1.25
$ws.Range("X2").Value = 12
$ws.Range("AO2").Value = 9
$ws.Range("BD2").Value = 150

# Row 3
$ws.Range("G3").Value = 3.5
$ws.Range("H3").Value = 3.25
$ws.Range("I3").Value = 2.15
$ws.Range("L3").Value = 2.88
$ws.Range("O3").Value = 1.33
$ws.Range("P3").Value = 3.4
$ws.Range("Q3").Value = 2.08
$ws.Range("R3").Value = 1.82
$ws.Range("S3").Value = 1.44
$ws.Range("T3").Value = 2.63
$ws.Range("AC3").Value = 9.5
$ws.Range("AD3").Value = 6
$ws.Range("AE3").Value = 13
$ws.Range("AM3").Value = 29
$ws.Range("AT3").Value = 2.63
$ws.Range("AZ3").Value = 23

# Row 5
$ws.Range("H5").Value = 3.2
$ws.Range("I5").Value = 4.2
$ws.Range("J5").Value = 2.75
$ws.Range("K5").Value = 1.95
$ws.Range("L5").Value = 4.75
$ws.Range("M5").Value = 1.08
$ws.Range("N5").Value = 8
$ws.Range("O5").Value = 1.44
$ws.Range("P5").Value = 2.63
$ws.Range("Q5").Value = 2.4
$ws.Range("R5").Value = 1.53
$ws.Range("S5").Value = 1.53
$ws.Range("T5").Value = 2.38
$ws.Range("U5").Value = 2.1
$ws.Range("V5").Value = 1.67
$ws.Range("W5").Value = 6
$ws.Range("AA5").Value = 19
$ws.Range("AB5").Value = 34
$ws.Range("AC5").Value = 7
$ws.Range("AD5").Value = 6
$ws.Range("AE5").Value = 19
$ws.Range("AF5").Value = 67
$ws.Range("AG5").Value = 501
$ws.Range("AH5").Value = 9
$ws.Range("AO5").Value = 11
$ws.Range("AP5").Value = 26
$ws.Range("AR5").Value = 67
$ws.Range("AS5").Value = 251
$ws.Range("AT5").Value = 2.38
$ws.Range("AU5").Value = 9
$ws.Range("AV5").Value = 67
$ws.Range("AY5").Value = 23
$ws.Range("BA5").Value = 81
$ws.Range("BB5").Value = 126
$ws.Range("BC5").Value = 351
$ws.Range("BD5").Value = 151

# Row 7
$ws.Range("O7").Value = 1.1
$ws.Range("P7").Value = 7
$ws.Range("U7").Value = 2
$ws.Range("V7").Value = 1.73
$ws.Range("Y7").Value = 11
$ws.Range("AK7").Value = 251
$ws.Range("AO7").Value = 4.75
$ws.Range("AU7").Value = 11
$ws.Range("BA7").Value = 301
$ws.Range("BB7").Value = 251
$ws.Range("BC7").Value = 351

# Row 9
$ws.Range("G9").Value = 2
$ws.Range("H9").Value = 3.25
$ws.Range("I9").Value = 3.6
$ws.Range("J9").Value = 2.75
$ws.Range("K9").Value = 2.2
$ws.Range("M9").Value = 1.05
$ws.Range("N9").Value = 11
$ws.Range("O9").Value = 1.25
$ws.Range("P9").Value = 3.75
$ws.Range("Q9").Value = 1.83
$ws.Range("R9").Value = 2.03
$ws.Range("W9").Value = 8
$ws.Range("AB9").Value = 26
$ws.Range("AC9").Value = 10
$ws.Range("AD9").Value = 6.5
$ws.Range("AI9").Value = 19
$ws.Range("AO9").Value = 11
$ws.Range("AZ9").Value = 26

# Row 18
$ws.Range("L18").Value = 13
$ws.Range("Q18").Value = 1.48
$ws.Range("R18").Value = 2.6
$ws.Range("AF18").Value = 81
$ws.Range("AK18").Value = 301
$ws.Range("AO18").Value = 4.75
$ws.Range("AR18").Value = 34
$ws.Range("AY18").Value = 67

# Row 19
$ws.Range("G19").Value = 4.5
$ws.Range("H19").Value = 3.6
$ws.Range("I19").Value = 1.8
$ws.Range("J19").Value = 5
$ws.Range("K19").Value = 2.1
$ws.Range("L19").Value = 2.4
$ws.Range("Q19").Value = 2.05
$ws.Range("R19").Value = 1.75
$ws.Range("W19").Value = 12
$ws.Range("X19").Value = 23
$ws.Range("Z19").Value = 51
$ws.Range("AA19").Value = 41
$ws.Range("AC19").Value = 9.5
$ws.Range("AD19").Value = 7
$ws.Range("AI19").Value = 8
$ws.Range("AJ19").Value = 8.5
$ws.Range("AL19").Value = 15
$ws.Range("AO19").Value = 26
$ws.Range("AR19").Value = 126
$ws.Range("AV19").Value = 51
$ws.Range("AY19").Value = 9.5
$ws.Range("AZ19").Value = 21
$ws.Range("BA19").Value = 34
$ws.Range("BC19").Value = 151
